$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 127, pushing the existing
# rows 127-137 down to 129-139 (carrying their formatting with them).
$ws.Rows("127:128").Insert()

# Row 127 (new): Primera, $/caja 70 unidades
$ws.Range("A127").Value = 2
$ws.Range("B127").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C127").Value = "Coquimbo"
$ws.Range("D127").Value = 44615
$ws.Range("D127").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E127").Value = 4
$ws.Range("F127").Value = 100112043
$ws.Range("G127").Value = "Pepino ensalada"
$ws.Range("H127").Value = "Sin especificar"
$ws.Range("I127").Value = "Primera"
$ws.Range("J127").Value = 500
$ws.Range("K127").Value = 12000
$ws.Range("L127").Value = 13000
$ws.Range("M127").Value = 12500
$ws.Range("N127").Value = "$/caja 70 unidades"
$ws.Range("O127").Value = "Provincia de Limarí"
$ws.Range("P127").Value = 179
$ws.Range("Q127").Value = 70
$ws.Range("R127").Value = "Hortaliza"

# Row 128 (new): Segunda, $/caja 100 unidades
$ws.Range("A128").Value = 2
$ws.Range("B128").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C128").Value = "Coquimbo"
$ws.Range("D128").Value = 44615
$ws.Range("D128").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E128").Value = 4
$ws.Range("F128").Value = 100112043
$ws.Range("G128").Value = "Pepino ensalada"
$ws.Range("H128").Value = "Sin especificar"
$ws.Range("I128").Value = "Segunda"
$ws.Range("J128").Value = 300
$ws.Range("K128").Value = 9000
$ws.Range("L128").Value = 10000
$ws.Range("M128").Value = 9500
$ws.Range("N128").Value = "$/caja 100 unidades"
$ws.Range("O128").Value = "Provincia de Limarí"
$ws.Range("P128").Value = 95
$ws.Range("Q128").Value = 100
$ws.Range("R128").Value = "Hortaliza"
